# Update gh-pages to output generated at 456a3b4
# Applies updated "want-to-go" headcount (column F) values across the
# 展览 (sheet1), 演出 (sheet2) and 全部类型 (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

# --- 展览 sheet -------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

$expoChanges = @{
    2  = 10200
    3  = 426
    5  = 22
    6  = 283
    7  = 187
    9  = 767
    11 = 1224
    12 = 1059
    13 = 3174
    14 = 2373
    16 = 2102
    17 = 2102
    21 = 1585
    22 = 565
    23 = 57
    24 = 243
    25 = 8
    27 = 237
    28 = 48
    31 = 370
    32 = 586
    33 = 50
    34 = 242
    35 = 5
    37 = 168
    38 = 393
    39 = 1683
    40 = 124
    41 = 428
    43 = 444
    44 = 971
    46 = 354
}

foreach ($row in $expoChanges.Keys) {
    $wsExpo.Range("F$row").Value = $expoChanges[$row]
}

# --- 演出 sheet -------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")

$showChanges = @{
    4 = 39
    9 = 7
}

foreach ($row in $showChanges.Keys) {
    $wsShow.Range("F$row").Value = $showChanges[$row]
}

# --- 全部类型 sheet ----------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$allChanges = @{
    2  = 10200
    3  = 426
    6  = 22
    8  = 283
    9  = 187
    11 = 767
    12 = 1059
    13 = 3174
    14 = 2373
    15 = 2102
    16 = 2102
    20 = 1585
    21 = 565
    22 = 57
    23 = 243
    24 = 8
    26 = 237
    27 = 48
    30 = 370
    31 = 586
    32 = 39
    35 = 50
    36 = 242
    37 = 5
    39 = 169
    41 = 393
    42 = 1683
    43 = 124
    45 = 428
    47 = 444
    48 = 971
    49 = 354
    50 = 7
}

foreach ($row in $allChanges.Keys) {
    $wsAll.Range("F$row").Value = $allChanges[$row]
}
